# spring 23 inputs complete
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New matchup rows (A:B:C:D) being appended after the existing data.
$data = @(
    @(6,3,4,0),
    @(7,3,2,0),
    @(3,2,5,0),
    @(2,1,4,2),
    @(5,0,7,2),
    @(3,1,4,2),
    @(6,0,6,2),
    @(4,0,5,2),
    @(4,3,3,0),
    @(4,1,4,2),
    @(5,0,5,2),
    @(4,2,4,0),
    @(4,0,5,2),
    @(5,0,5,2),
    @(3,1,5,2),
    @(4,0,4,3),
    @(3,0,5,3),
    @(5,0,6,2),
    @(3,2,4,1),
    @(5,0,5,2),
    @(3,1,3,2),
    @(7,0,6,2)
)

$startRow = 1900
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Range("A" + $r).Value = $row[0]
    $ws.Range("B" + $r).Value = $row[1]
    $ws.Range("C" + $r).Value = $row[2]
    $ws.Range("D" + $r).Value = $row[3]
}

$lastRow = $startRow + $data.Count - 1
$nextCell = "A" + ($lastRow + 1)

# Scroll/select so the view lands just past the newly entered rows
# (mirrors the author scrolling down after finishing data entry).
$win = $excel.ActiveWindow
$win.ScrollRow = $startRow - 4
$ws.Range($nextCell).Select()
